$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.622.44'
$ws.Range("E2").Value = '  -1.01%  '
$ws.Range("D3").Value = '1.658.40'
$ws.Range("E3").Value = '  -2.92%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.24'
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9979'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3652'
$ws.Range("E7").Value = '  -2.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.79'
$ws.Range("E8").Value = '  -5.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3272'
$ws.Range("E9").Value = '  -4.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.132'
$ws.Range("E10").Value = '  -7.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07061'
$ws.Range("E11").Value = '  -6.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9968'
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.997'
$ws.Range("E13").Value = '  -5.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.54'
$ws.Range("E14").Value = '  -8.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.625'
$ws.Range("E15").Value = '  -6.79%  '
$ws.Range("D16").Value = '1.654.92'
$ws.Range("E16").Value = '  -3.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001047'
$ws.Range("E17").Value = '  -7.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06630'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9972'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '79.15'
$ws.Range("E20").Value = '  -5.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.962'
$ws.Range("E21").Value = '  -6.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.83'
$ws.Range("E22").Value = '  -8.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.64'
$ws.Range("E23").Value = '  -3.68%  '
$ws.Range("D24").Value = '24.585.75'
$ws.Range("E24").Value = '  -1.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.466'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.385'
$ws.Range("E26").Value = '  -15.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.04'
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.69'
$ws.Range("E28").Value = '  -8.43%  '
$ws.Range("E29").Value = '  -2.01%  '
$ws.Range("D30").Value = '1.839.49'
$ws.Range("E30").Value = '  -3.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '124.99'
$ws.Range("E31").Value = '  -6.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.064'
$ws.Range("E32").Value = '  -3.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.829'
$ws.Range("E33").Value = '  -14.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08488'
$ws.Range("E34").Value = '  -3.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.681'
$ws.Range("E35").Value = '  -4.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.33'
$ws.Range("E36").Value = '  -11.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.241'
$ws.Range("E37").Value = '  -6.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.267'
$ws.Range("E38").Value = '  -0.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06048'
$ws.Range("E39").Value = '  -9.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02239'
$ws.Range("E40").Value = '  -7.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2076'
$ws.Range("E41").Value = '  -7.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.167'
$ws.Range("E42").Value = '  -11.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9976'
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5931'
$ws.Range("E44").Value = '  -8.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.866'
$ws.Range("E45").Value = '  +0.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.72'
$ws.Range("E46").Value = '  -8.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5640'
$ws.Range("E47").Value = '  -8.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.50'
$ws.Range("E48").Value = '  -4.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.960'
$ws.Range("E49").Value = '  -8.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06971'
$ws.Range("E50").Value = '  -4.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.194'
$ws.Range("E51").Value = '  -3.16%  '
